# Append new log rows to the worksheet (thread safe log textbox data dump)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(42601.977083333331, "Bag", 25, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0),
    @(42601.988379629627, "Bag", 20, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0),
    @(42601.988946759258, "Bag", 19, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0),
    @(42601.991145833330, "Bag", 19, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0),
    @(42601.995023148149, "Bag", 66, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0),
    @(42601.997997685183, "Bag", 15, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0),
    @(42601.999166666668, "Bag", 72, 17, 2, 0, 1, 0, 100, 2, 0, 100, 0)
)

$startRow = 13
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
